$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Swap the match data (columns F:V) between row 117 and row 118 ---
# Row 117 currently holds the Lecco-Spezia match result/odds; row 118 holds the
# Palermo-Brescia match result/odds. After the edit, row 117 should hold the
# Palermo-Brescia data and row 118 should hold the Lecco-Spezia data. Columns
# A:E (Indice, pais, torneio, temporada, data_partida) are already identical
# between the two rows, so only F:V need to be exchanged.
$row117 = $ws.Range("F117:V117")
$row118 = $ws.Range("F118:V118")

$buffer117 = $row117.Value()
$buffer118 = $row118.Value()

$row117.Value = $buffer118
$row118.Value = $buffer117

# --- Step 2: Append the new match (Venezia vs Catanzaro) as row 119 ---
# Copy the formatting used by the preceding rows for the index column (bold,
# centered, bordered) and the match-date column (custom date/time format).
$ws.Range("A118").Copy()
$ws.Range("A119").PasteSpecial(-4122)
$ws.Range("E118").Copy()
$ws.Range("E119").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A119").Value = 118
$ws.Range("B119").Value = "italy"
$ws.Range("C119").Value = "serie-b"
$ws.Range("D119").Value = "2023-2024"
$ws.Range("E119").Value = 45240.85416666666
$ws.Range("F119").Value = "Venezia"
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = "Catanzaro"
$ws.Range("I119").Value = 1
$ws.Range("J119").Value = 2.02
$ws.Range("K119").Value = "04/11/2023 14:12"
$ws.Range("L119").Value = 2.02
$ws.Range("M119").Value = "10/11/2023 20:01"
$ws.Range("N119").Value = 3.58
$ws.Range("O119").Value = "04/11/2023 14:12"
$ws.Range("P119").Value = 3.53
$ws.Range("Q119").Value = "10/11/2023 20:01"
$ws.Range("R119").Value = 3.88
$ws.Range("S119").Value = "04/11/2023 14:12"
$ws.Range("T119").Value = 4
$ws.Range("U119").Value = "10/11/2023 20:01"
$ws.Range("V119").Value = "https://www.betexplorer.com/football/italy/serie-b/venezia-catanzaro/p84Xyfrk/"
